$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.780.90"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "3.165.30"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'615.68"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").Value = "'145.99"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.163.65"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("D9").Value = "'0.531"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "'0.153"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "'5.51"
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D12").Value = "'0.475"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").Value = "'0.0000261"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "'35.89"
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").Value = "3.685.38"
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("D17").Value = "64.749.15"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").Value = "3.165.70"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").Value = "'6.93"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "'480.08"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "'14.74"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("D23").Value = "'7.98"
$ws.Range("E23").Value = "  +2.68%  "
$ws.Range("D24").Value = "'13.80"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'84.72"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -3.21%  "
$ws.Range("D28").Value = "'8.65"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("E29").Value = "  -6.53%  "
$ws.Range("D30").Value = "'6.91"
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("D31").Value = "'2.10"
$ws.Range("E31").Value = "  -6.56%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "'2.71"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("D34").Value = "'26.70"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  +3.46%  "
$ws.Range("D36").Value = "0.0₃0792"
$ws.Range("E36").Value = "  +5.95%  "
$ws.Range("D37").Value = "'6.04"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").Value = "'53.20"
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("D40").Value = "'466.41"
$ws.Range("E40").Value = "  +3.60%  "
$ws.Range("D41").Value = "'0.0401"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").Value = "'0.121"
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("D43").Value = "'8.41"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("D44").Value = "2.844.93"
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("D45").Value = "'2.34"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "'2.45"
$ws.Range("E47").Value = "  +5.92%  "
$ws.Range("D48").Value = "'26.76"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "'34.45"
$ws.Range("E51").Value = "  +4.39%  "
